$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "singleton" (single-cell, no-op) merged ranges left over from the
# old ClosedXML version. Excel/ClosedXML no longer emits 1x1 merges.
$ws.Range("C2:C2").UnMerge()
$ws.Range("F2:F2").UnMerge()
$ws.Range("B4:B4").UnMerge()
$ws.Range("D4:D4").UnMerge()
$ws.Range("E4:E4").UnMerge()
$ws.Range("B5:B5").UnMerge()
$ws.Range("E5:E5").UnMerge()
$ws.Range("F5:F5").UnMerge()
$ws.Range("B6:B6").UnMerge()
$ws.Range("E6:E6").UnMerge()
$ws.Range("B7:B7").UnMerge()
$ws.Range("D7:D7").UnMerge()
$ws.Range("E7:E7").UnMerge()
$ws.Range("B8:B8").UnMerge()

# B5 and B6 become a genuine two-cell merge.
$ws.Range("B5:B6").Merge()

# F4 carried a stray shared-string type with no backing value (an artifact of
# the old ClosedXML writer). Round-trip the F3:F4 merge so the cell is
# rewritten as a truly empty cell (no value, no type).
$ws.Range("F3:F4").UnMerge()
$ws.Range("F4").ClearContents()
$ws.Range("F3:F4").Merge()

# F8:F9 was merged around an (empty) cell. Now F9 receives its own value, so
# the merge between F8 and F9 is removed.
$ws.Range("F8:F9").UnMerge()

# Give F9 the default ("General") style instead of the date number format
# (style index 1) it inherited from F8, then assign its new text value.
$ws.Range("A1").Copy()
$ws.Range("F9").PasteSpecial(-4122)
$ws.Range("F9").Value = "21.02.2018 0:00:00"
